$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Screen 1 (B2/G2 pair) ---
$ws.Range("B2").Value = 240
$ws.Range("G2").Value = 360

# --- Screen 2 (B3/G3 pair) ---
$ws.Range("B3").Value = 427
$ws.Range("G3").Value = 640

# --- New reference row: device name + its resolution + computed ratio ---
$ws.Range("I4").Value = "Sunmi V1s"
$ws.Range("J4").Value = 360
$ws.Range("K4").Value = 640
$ws.Range("L4").Value = 4.1225529999999999

# --- Diagonal size used for PPI computation ---
$ws.Range("C6").Value = 2.75

# --- Page setup (paper size 256 = custom/user-defined, portrait) ---
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1

# --- Selection moves to L4 after the edits ---
[void]$ws.Range("L4").Select()
